$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 48) to the daily log sheet, extending the
# table from A1:D47 to A1:D48 with one more timestamped ranking entry.
#
# Column A holds a date-like string ("2025/10/02") that must stay a plain
# text value (as all the other rows in the sheet are), not get auto-converted
# into a date serial number by Excel's smart input parsing. Temporarily
# forcing a text number format while assigning the value, then clearing the
# formatting back to the sheet's default (unstyled) look, achieves that.
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "2025/10/02"
$ws.Range("A48").ClearFormats()

$ws.Range("B48").Value = "木"
$ws.Range("C48").Value = 7
$ws.Range("D48").Value = 3
